$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 164, shifting rows 164:170 down to 165:171.
$ws.Rows(164).Insert()

# Populate the newly inserted row 164 with the new weekly record.
$ws.Range("A164").Value = 4
$ws.Range("B164").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C164").Value = "Los Lagos"
$ws.Range("D164").Value = 44509
$ws.Range("E164").Value = 10
$ws.Range("F164").Value = 100112003
$ws.Range("G164").Value = "Ajo"
$ws.Range("H164").Value = "Chino"
$ws.Range("I164").Value = "Primera"
$ws.Range("J164").Value = 240
$ws.Range("K164").Value = 21000
$ws.Range("L164").Value = 22000
$ws.Range("M164").Value = 21500
$ws.Range("N164").Value = "`$/caja 10 kilos"
$ws.Range("O164").Value = "China"
$ws.Range("P164").Value = 2150
$ws.Range("Q164").Value = 10
$ws.Range("R164").Value = "Hortaliza"
